{"js": "// Replace the date line and every \"NNN\u00f7N=\" division prompt in the table\n// with the new values from the day's worksheet, per the commit diff.\nconst replacements = [\n  [\"2024-06-23 Sunday\", \"2024-06-24 Monday\"],\n  [\"748\u00f73=\", \"616\u00f73=\"],\n  [\"486\u00f79=\", \"832\u00f78=\"],\n  [\"215\u00f73=\", \"165\u00f72=\"],\n  [\"766\u00f76=\", \"406\u00f73=\"],\n  [\"600\u00f76=\", \"852\u00f78=\"],\n  [\"573\u00f79=\", \"862\u00f73=\"],\n  [\"309\u00f74=\", \"137\u00f78=\"],\n  [\"344\u00f79=\", \"498\u00f73=\"],\n  [\"649\u00f74=\", \"853\u00f79=\"],\n  [\"501\u00f78=\", \"206\u00f76=\"],\n  [\"397\u00f72=\", \"830\u00f75=\"],\n  [\"487\u00f78=\", \"606\u00f74=\"],\n  [\"411\u00f78=\", \"841\u00f73=\"],\n  [\"415\u00f75=\", \"725\u00f77=\"],\n  [\"195\u00f74=\", \"993\u00f75=\"],\n  [\"881\u00f73=\", \"631\u00f77=\"],\n  [\"786\u00f75=\", \"976\u00f76=\"],\n  [\"510\u00f76=\", \"411\u00f79=\"],\n  [\"975\u00f75=\", \"485\u00f74=\"],\n  [\"285\u00f74=\", \"331\u00f72=\"],\n  [\"137\u00f72=\", \"627\u00f75=\"],\n  [\"121\u00f72=\", \"572\u00f73=\"],\n  [\"586\u00f72=\", \"107\u00f72=\"],\n  [\"557\u00f72=\", \"901\u00f73=\"],\n  [\"283\u00f75=\", \"507\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"NNN\u00f7N=\" division prompt to the\n# next day's values, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-06-23 Sunday', '2024-06-24 Monday'),\n    @('748\u00f73=', '616\u00f73='),\n    @('486\u00f79=', '832\u00f78='),\n    @('215\u00f73=', '165\u00f72='),\n    @('766\u00f76=', '406\u00f73='),\n    @('600\u00f76=', '852\u00f78='),\n    @('573\u00f79=', '862\u00f73='),\n    @('309\u00f74=', '137\u00f78='),\n    @('344\u00f79=', '498\u00f73='),\n    @('649\u00f74=', '853\u00f79='),\n    @('501\u00f78=', '206\u00f76='),\n    @('397\u00f72=', '830\u00f75='),\n    @('487\u00f78=', '606\u00f74='),\n    @('411\u00f78=', '841\u00f73='),\n    @('415\u00f75=', '725\u00f77='),\n    @('195\u00f74=', '993\u00f75='),\n    @('881\u00f73=', '631\u00f77='),\n    @('786\u00f75=', '976\u00f76='),\n    @('510\u00f76=', '411\u00f79='),\n    @('975\u00f75=', '485\u00f74='),\n    @('285\u00f74=', '331\u00f72='),\n    @('137\u00f72=', '627\u00f75='),\n    @('121\u00f72=', '572\u00f73='),\n    @('586\u00f72=', '107\u00f72='),\n    @('557\u00f72=', '901\u00f73='),\n    @('283\u00f75=', '507\u00f73=')\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
